$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '257.71'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '5.12%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '28.01'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-2.91%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.214'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-1.28%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05903'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '3.34%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.703'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '1.26%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.8691'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '1.85%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.040'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '21.20%'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1412'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '3.05%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07187'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '2.16%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.03150'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-0.27%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.09213'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-0.96%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.001538'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '0.37%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0006077'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-93.95%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.005911'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-2.27%'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.229'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '1.53%'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.226'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '2.36%'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3121'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-1.37%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.03633'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '8.62%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1290'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '1.01%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.526'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '0.67%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04205'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '2.84%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.1365'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-1.04%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.001218'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-0.54%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.004548'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '9.77%'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0001199'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '-0.10%'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0001472'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '1.52%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03835'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '1.89%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.005417'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '4.51%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1103'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '3.67%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002299'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-6.16%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.01066'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '14.06%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005428'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000750'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-0.11%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.08546'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '13.89%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.002136'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '-12.51%'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002099'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-0.11%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0001999'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.11%'
